$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.249.30'
$ws.Range('E2').Value = '  -1.62%  '

$ws.Range('D3').Value = '2.921.00'
$ws.Range('E3').Value = '  -2.34%  '

$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '373.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.40%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.65%  '

$ws.Range('E7').Value = '  -3.62%  '

$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.591'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.08'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.99%  '

$ws.Range('E11').Value = '  +0.43%  '

$ws.Range('E12').Value = '  -2.28%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.18%  '

$ws.Range('D14').Value = '3.378.34'
$ws.Range('E14').Value = '  -2.44%  '

$ws.Range('E15').Value = '  -3.70%  '

$ws.Range('D16').Value = '2.920.37'
$ws.Range('E16').Value = '  -2.01%  '

$ws.Range('E17').Value = '  -9.09%  '

$ws.Range('D18').Value = '51.193.62'
$ws.Range('E18').Value = '  -1.84%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.63%  '

$ws.Range('D22').Value = '0.0₃0947'
$ws.Range('E22').Value = '  -2.39%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.40%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '260.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.10%  '

$ws.Range('E26').Value = '  -4.86%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.99%  '

$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.50%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.95%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.102'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.82%  '

$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.82%  '

$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.52%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.64'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.11%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.61%  '

$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.40%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0425'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.25%  '

$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.50%  '

$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.35%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.22%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.34%  '

$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.113'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.29%  '

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.17%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.74%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.025.39'
$ws.Range('E46').Value = '  -4.67%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.48%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.89%  '

$ws.Range('B49').Value = 'TheGraph'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.246'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.96%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '3.212.13'
$ws.Range('E50').Value = '  -2.26%  '

$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0316'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.98%  '
